# BLSA and BGISA debugging
# Add dark mode and tooltip improvements
#
# Re-classifies the "Enforceability" column (F) on Sheet1: most rows that
# were marked "Unknown" are now resolved to a concrete rating, and the
# "Weak" bucket gains a couple of more specific siblings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> new Enforceability (column F) value.
$enforceability = @{
    5  = "Strong"
    25 = "Strong"
    26 = "Strong"
    27 = "Strong"
    28 = "Strong"
    29 = "Strong"
    30 = "Strong"
    31 = "Strong"
    32 = "Strong"
    33 = "Strong"
    34 = "Strong"
    35 = "Strong"
    36 = "Strong"
    37 = "Strong"
    38 = "Strong"
    39 = "Weak"
    40 = "Strong"
    41 = "Strong"
    42 = "Strong"
    43 = "Strong"
    44 = "Strong"
    45 = "Strong"
    46 = "Strong"
    47 = "Weak (currently)"
    48 = "Strong"
    49 = "Not Applicable"
    50 = "Not Applicable"
    51 = "Not Applicable"
    52 = "Not Applicable"
}

foreach ($row in $enforceability.Keys) {
    $ws.Cells.Item($row, 6).Value = $enforceability[$row]
}

# Restore the scroll position / selection the author left the sheet in.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F49").Select() | Out-Null

# Window placement for the workbook itself.
$excel.ActiveWindow.Left = -250
$excel.ActiveWindow.Top = 70
$excel.ActiveWindow.Width = 9710
$excel.ActiveWindow.Height = 11280
